$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: fill in the new log entry ---

# Full rich text for A12, built from 6 runs with different formatting.
$full = "Profile Pagina" + ": Klant en klusjesman krijgen verschillende pagina te zien waarop andere gegevens te vinden zijn" + " [maincontroller]" + ". Ook nieuwe beta querries toegevoegd om deze verschillende data op te halen " + "[taskrepo]" + "."

$cellA12 = $ws.Range("A12")
$cellA12.Value = $full

# Run 1 (bold): "Profile Pagina" -> chars 1-14
$cellA12.Characters(1, 14).Font.Bold = $true

# Run 2 (regular): ": Klant en klusjesman ... te vinden zijn" -> chars 15-110
$cellA12.Characters(15, 96).Font.Bold = $false

# Run 3 (italic): " [maincontroller]" -> chars 111-127
$cellA12.Characters(111, 17).Font.Italic = $true

# Run 4 (regular): ". Ook nieuwe beta querries ... op te halen " -> chars 128-204
$cellA12.Characters(128, 77).Font.Italic = $false

# Run 5 (italic): "[taskrepo]" -> chars 205-214
$cellA12.Characters(205, 10).Font.Italic = $true

# Run 6 (regular): "." -> char 215
$cellA12.Characters(215, 1).Font.Italic = $false

# Date worked (2024-12-12, Excel serial date number 45638)
$ws.Range("B12").Value = 45638

# Hours worked
$ws.Range("C12").Value = 0.5

# Row height (auto-expanded by Excel because of the wrapped multi-run text)
$ws.Rows(12).RowHeight = 109.2

# --- Update the active selection to D12, matching the end-user's last click ---
$ws.Range("D12").Select()
